$wb = $excel.ActiveWorkbook

# ==================================================================
# Helper: force a value to be written as literal TEXT (preserving
# leading zeros / trailing decimal zeros that Excel's implicit number
# coercion would otherwise eat) without leaving any stray cell style
# behind. Routes the text through a scratch cell as a string formula,
# then pastes-special just the resulting value into the destination.
# ==================================================================
$scratchSheet = $wb.Worksheets.Item(1)
$scratch = $scratchSheet.Range("ZZ1")

function Set-TextValue {
    param($range, [string]$text)
    $esc = $text -replace '"', '""'
    $scratch.Formula = '="' + $esc + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

# ------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new 2022-Q3 row at the top of the
#    data (row 2), pushing the existing quarters down by one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make room: give the new row 9 the same style as the existing A-column
# data cells by copying A8 -> A9 (also value, overwritten right after).
$summary.Range("A8").Copy($summary.Range("A9"))

# Rewrite column A (index numbers 0..7) for all 8 data rows.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# Row 2: new 2022-Q3 entry
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 29
$summary.Range("D2").Value = 10.31

# Row 3: was row2 (2022-Q2)
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 23
$summary.Range("D3").Value = 11.8

# Row 4: was row3 (2022-Q1)
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 16
$summary.Range("D4").Value = 10.76

# Row 5: was row4 (2021-Q4)
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 23
$summary.Range("D5").Value = 12.69

# Row 6: was row5 (2021-Q3)
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 6
$summary.Range("D6").Value = 3.16

# Row 7: was row6 (2021-Q2)
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 3
$summary.Range("D7").Value = 2.01

# Row 8: was row7 (2021-Q1)
$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 3
$summary.Range("D8").Value = 2.03

# Row 9 (new): was row8 (2020-Q4)
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 6
$summary.Range("D9").Value = 2.3

# ------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" sheet right after "总计" (i.e. before
#    the current "2022-Q2" sheet), containing the per-fund breakdown.
# ------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($oldQ2)
$q3.Name = "2022-Q3"

# Header row, copying style from the old 2022-Q2 sheet's header cells so
# formatting (bold/border/center) matches the rest of the workbook.
$oldQ2.Range("B1:H1").Copy($q3.Range("B1:H1"))
Set-TextValue $q3.Range("B1") "基金代码"
Set-TextValue $q3.Range("C1") "基金名称"
Set-TextValue $q3.Range("D1") "基金规模"
Set-TextValue $q3.Range("E1") "股票总仓位"
Set-TextValue $q3.Range("F1") "仓位占比"
Set-TextValue $q3.Range("G1") "持有市值(亿元)"
Set-TextValue $q3.Range("H1") "仓位排名"

function Set-FundRow {
    param($ws, $row, $a, $b, $c, $d, $e, $f, $g, $h)
    $ws.Range("A$row").Value = $a
    Set-TextValue $ws.Range("B$row") $b
    Set-TextValue $ws.Range("C$row") $c
    Set-TextValue $ws.Range("D$row") $d
    Set-TextValue $ws.Range("E$row") $e
    Set-TextValue $ws.Range("F$row") $f
    Set-TextValue $ws.Range("G$row") $g
    $ws.Range("H$row").Value = $h
}

# Apply the A-column ("s=2") style used throughout the workbook by
# copying it from the old sheet's A2 cell across the new data rows.
for ($r = 2; $r -le 30; $r++) {
    $oldQ2.Range("A2").Copy($q3.Range("A$r"))
}

Set-FundRow $q3 2  0  "320003" "诺安先锋混合A" "40.90" "76.11" "7.04" "2.8794" 1
Set-FundRow $q3 3  1  "000362" "国泰聚信价值优势灵活配置混合A" "27.52" "89.04" "5.72" "1.5741" 3
Set-FundRow $q3 4  2  "020010" "国泰金牛创新混合" "13.26" "86.27" "6.09" "0.8075" 1
Set-FundRow $q3 5  3  "160106" "南方高增长混合（LOF）" "15.29" "91.51" "5.12" "0.7828" 7
Set-FundRow $q3 6  4  "000363" "国泰聚信价值优势灵活配置混合C" "13.07" "89.04" "5.72" "0.7476" 3
Set-FundRow $q3 7  5  "001743" "诺安优选回报灵活配置混合" "13.65" "73.02" "4.51" "0.6156" 1
Set-FundRow $q3 8  6  "008415" "国泰大制造两年持有期混合" "10.30" "90.83" "4.83" "0.4975" 4
Set-FundRow $q3 9  7  "012173" "国泰兴泽优选一年持有期混合A" "8.41" "88.23" "5.22" "0.4390" 2
Set-FundRow $q3 10 8  "007835" "国泰鑫睿混合" "8.30" "79.49" "4.86" "0.4034" 1
Set-FundRow $q3 11 9  "012174" "国泰兴泽优选一年持有期混合C" "6.17" "88.23" "5.22" "0.3221" 2
Set-FundRow $q3 12 10 "160105" "南方积极配置混合（LOF）" "5.42" "91.29" "5.21" "0.2824" 7
Set-FundRow $q3 13 11 "005244" "国泰聚优价值灵活配置混合A" "4.61" "87.30" "4.61" "0.2125" 3
Set-FundRow $q3 14 12 "003131" "国寿安保强国智造灵活配置混合" "4.44" "92.60" "2.45" "0.1088" 5
Set-FundRow $q3 15 13 "012442" "永赢稳健增长一年持有期混合E" "9.56" "26.04" "1.09" "0.1042" 7
Set-FundRow $q3 16 14 "008185" "诺安研究优选混合A" "1.36" "93.89" "7.43" "0.1010" 2
Set-FundRow $q3 17 15 "005245" "国泰聚优价值灵活配置混合C" "1.80" "87.30" "4.61" "0.0830" 3
Set-FundRow $q3 18 16 "012621" "诺安先锋混合C" "1.13" "76.11" "7.04" "0.0796" 1
Set-FundRow $q3 19 17 "000554" "南方中国梦灵活配置混合" "1.33" "90.51" "5.23" "0.0696" 7
Set-FundRow $q3 20 18 "009932" "永赢稳健增长一年持有期混合A" "4.93" "26.04" "1.09" "0.0537" 7
Set-FundRow $q3 21 19 "005683" "国寿安保华兴灵活配置混合" "2.00" "92.53" "2.17" "0.0434" 5
Set-FundRow $q3 22 20 "010797" "长城优选回报六个月持有期混合A" "2.83" "28.94" "1.17" "0.0331" 6
Set-FundRow $q3 23 21 "014321" "德邦周期精选混合A" "0.62" "92.98" "3.32" "0.0206" 9
Set-FundRow $q3 24 22 "010857" "宝盈祥乐一年持有期混合型证券投资基金A" "1.08" "27.77" "1.87" "0.0202" 3
Set-FundRow $q3 25 23 "008324" "宝盈祥利稳健配置混合A" "0.55" "27.49" "1.78" "0.0098" 3
Set-FundRow $q3 26 24 "014497" "诺安研究优选混合C" "0.08" "93.89" "7.43" "0.0059" 2
Set-FundRow $q3 27 25 "008325" "宝盈祥利稳健配置混合C" "0.31" "27.49" "1.78" "0.0055" 3
Set-FundRow $q3 28 26 "010798" "长城优选回报六个月持有期混合C" "0.30" "28.94" "1.17" "0.0035" 6
Set-FundRow $q3 29 27 "010858" "宝盈祥乐一年持有期混合型证券投资基金C" "0.06" "27.77" "1.87" "0.0011" 3
Set-FundRow $q3 30 28 "014322" "德邦周期精选混合C" "0.03" "92.98" "3.32" "0.0010" 9

# Clean up the scratch cell used for text coercion.
$scratch.ClearContents()
